$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated loading-percent simulation results for the 380 kV case (rows 2-25)
$ws.Cells.Item(2, 3).Value2 = 3.763632882215974
$ws.Cells.Item(2, 4).Value2 = 8.260638386953504
$ws.Cells.Item(2, 5).Value2 = 14.46166996205858
$ws.Cells.Item(2, 6).Value2 = 30.28390284961943
$ws.Cells.Item(2, 7).Value2 = 35.39245372831265
$ws.Cells.Item(2, 8).Value2 = 15.64219040100228
$ws.Cells.Item(2, 9).Value2 = 22.88485295830392
$ws.Cells.Item(2, 10).Value2 = 9.925494871914859
$ws.Cells.Item(2, 13).Value2 = 59.58928092560715

$ws.Cells.Item(3, 3).Value2 = 3.773333365486835
$ws.Cells.Item(3, 4).Value2 = 8.141758435685299
$ws.Cells.Item(3, 5).Value2 = 14.06650873965341
$ws.Cells.Item(3, 6).Value2 = 30.90958678455896
$ws.Cells.Item(3, 7).Value2 = 35.82273597497176
$ws.Cells.Item(3, 8).Value2 = 15.84584562050025
$ws.Cells.Item(3, 9).Value2 = 23.15637632281554
$ws.Cells.Item(3, 10).Value2 = 9.808727376801206
$ws.Cells.Item(3, 13).Value2 = 56.22032196012083

$ws.Cells.Item(4, 3).Value2 = 3.780616860138308
$ws.Cells.Item(4, 4).Value2 = 8.068672830321583
$ws.Cells.Item(4, 5).Value2 = 13.82282467702637
$ws.Cells.Item(4, 6).Value2 = 31.31789444930932
$ws.Cells.Item(4, 7).Value2 = 36.13037342602779
$ws.Cells.Item(4, 8).Value2 = 15.97942013211202
$ws.Cells.Item(4, 9).Value2 = 23.34027986026106
$ws.Cells.Item(4, 10).Value2 = 9.740379806031601
$ws.Cells.Item(4, 13).Value2 = 54.03936855618517

$ws.Cells.Item(5, 3).Value2 = 3.783913943148021
$ws.Cells.Item(5, 4).Value2 = 8.038896161685358
$ws.Cells.Item(5, 5).Value2 = 13.72340242875046
$ws.Cells.Item(5, 6).Value2 = 31.49021488487936
$ws.Cells.Item(5, 7).Value2 = 36.26623968626305
$ws.Cells.Item(5, 8).Value2 = 16.03596248773102
$ws.Cells.Item(5, 9).Value2 = 23.4194300805347
$ws.Cells.Item(5, 10).Value2 = 9.713389867405112
$ws.Cells.Item(5, 13).Value2 = 53.12266142573471

$ws.Cells.Item(6, 3).Value2 = 3.784481154970941
$ws.Cells.Item(6, 4).Value2 = 8.033953025825241
$ws.Cells.Item(6, 5).Value2 = 13.70689067280514
$ws.Cells.Item(6, 6).Value2 = 31.5191830294474
$ws.Cells.Item(6, 7).Value2 = 36.28942223463248
$ws.Cells.Item(6, 8).Value2 = 16.0454777306154
$ws.Cells.Item(6, 9).Value2 = 23.43282381750947
$ws.Cells.Item(6, 10).Value2 = 9.708960854208124
$ws.Cells.Item(6, 13).Value2 = 52.96876570586585

$ws.Cells.Item(7, 3).Value2 = 3.780659999850255
$ws.Cells.Item(7, 4).Value2 = 8.068271186664367
$ws.Cells.Item(7, 5).Value2 = 13.82148411695674
$ws.Cells.Item(7, 6).Value2 = 31.32019457841497
$ws.Cells.Item(7, 7).Value2 = 36.13216379835021
$ws.Cells.Item(7, 8).Value2 = 15.98017418511939
$ws.Cells.Item(7, 9).Value2 = 23.3413304170381
$ws.Cells.Item(7, 10).Value2 = 9.740012293162156
$ws.Cells.Item(7, 13).Value2 = 54.02711821527924

$ws.Cells.Item(8, 3).Value2 = 3.766699200025426
$ws.Cells.Item(8, 4).Value2 = 8.219685465364085
$ws.Cells.Item(8, 5).Value2 = 14.32572193553682
$ws.Cells.Item(8, 6).Value2 = 30.49454324305893
$ws.Cells.Item(8, 7).Value2 = 35.53154247227013
$ws.Cells.Item(8, 8).Value2 = 15.71061872863652
$ws.Cells.Item(8, 9).Value2 = 22.97483738484329
$ws.Cells.Item(8, 10).Value2 = 9.884551105602256
$ws.Cells.Item(8, 13).Value2 = 58.45115253342351

$ws.Cells.Item(9, 3).Value2 = 3.750075874747136
$ws.Cells.Item(9, 4).Value2 = 8.514602174705242
$ws.Cells.Item(9, 5).Value2 = 15.30021913750341
$ws.Cells.Item(9, 6).Value2 = 29.07373096958377
$ws.Cells.Item(9, 7).Value2 = 34.71754707207322
$ws.Cells.Item(9, 8).Value2 = 15.25137087594184
$ws.Cells.Item(9, 9).Value2 = 22.39780556864919
$ws.Cells.Item(9, 10).Value2 = 10.19361497831874
$ws.Cells.Item(9, 13).Value2 = 66.22534478406004

$ws.Cells.Item(10, 3).Value2 = 3.744734586574623
$ws.Cells.Item(10, 4).Value2 = 8.728413118186094
$ws.Cells.Item(10, 5).Value2 = 16.0001994963063
$ws.Cells.Item(10, 6).Value2 = 28.16115452921745
$ws.Cells.Item(10, 7).Value2 = 34.36749735444975
$ws.Cells.Item(10, 8).Value2 = 14.95870039765151
$ws.Cells.Item(10, 9).Value2 = 22.06775261233429
$ws.Cells.Item(10, 10).Value2 = 10.43500021461604
$ws.Cells.Item(10, 13).Value2 = 71.37969834605813

$ws.Cells.Item(11, 3).Value2 = 3.743868004896
$ws.Cells.Item(11, 4).Value2 = 8.824731285816558
$ws.Cells.Item(11, 5).Value2 = 16.31384746199258
$ws.Cells.Item(11, 6).Value2 = 27.7771282010451
$ws.Cells.Item(11, 7).Value2 = 34.26762603219085
$ws.Cells.Item(11, 8).Value2 = 14.83585851282431
$ws.Cells.Item(11, 9).Value2 = 21.9396927577677
$ws.Cells.Item(11, 10).Value2 = 10.54763707272467
$ws.Cells.Item(11, 13).Value2 = 73.60293926670336

$ws.Cells.Item(12, 3).Value2 = 3.743770682705927
$ws.Cells.Item(12, 4).Value2 = 8.861042612508212
$ws.Cells.Item(12, 5).Value2 = 16.43183545192176
$ws.Cells.Item(12, 6).Value2 = 27.63642776952301
$ws.Cells.Item(12, 7).Value2 = 34.23878322194014
$ws.Cells.Item(12, 8).Value2 = 14.79087507763598
$ws.Cells.Item(12, 9).Value2 = 21.89452105407021
$ws.Cells.Item(12, 10).Value2 = 10.59067267408006
$ws.Cells.Item(12, 13).Value2 = 74.42734949238444

$ws.Cells.Item(13, 3).Value2 = 3.743781278324751
$ws.Cells.Item(13, 4).Value2 = 8.853229940548127
$ws.Cells.Item(13, 5).Value2 = 16.40646090699686
$ws.Cells.Item(13, 6).Value2 = 27.66651589710922
$ws.Cells.Item(13, 7).Value2 = 34.24458901626004
$ws.Cells.Item(13, 8).Value2 = 14.80049394317562
$ws.Cells.Item(13, 9).Value2 = 21.90409948523515
$ws.Cells.Item(13, 10).Value2 = 10.58138755974805
$ws.Cells.Item(13, 13).Value2 = 74.2505753950619

$ws.Cells.Item(14, 3).Value2 = 3.743855336709307
$ws.Cells.Item(14, 4).Value2 = 8.82772203875232
$ws.Cells.Item(14, 5).Value2 = 16.32357062686687
$ws.Cells.Item(14, 6).Value2 = 27.76545654681686
$ws.Cells.Item(14, 7).Value2 = 34.26507076754265
$ws.Cells.Item(14, 8).Value2 = 14.83212663285595
$ws.Cells.Item(14, 9).Value2 = 21.93590903089738
$ws.Cells.Item(14, 10).Value2 = 10.55117011029994
$ws.Cells.Item(14, 13).Value2 = 73.67111432492554

$ws.Cells.Item(15, 3).Value2 = 3.743930943204422
$ws.Cells.Item(15, 4).Value2 = 8.812075821270041
$ws.Cells.Item(15, 5).Value2 = 16.27269319132276
$ws.Cells.Item(15, 6).Value2 = 27.82668328219103
$ws.Cells.Item(15, 7).Value2 = 34.27879813985738
$ws.Cells.Item(15, 8).Value2 = 14.85170400817478
$ws.Cells.Item(15, 9).Value2 = 21.95583029165681
$ws.Cells.Item(15, 10).Value2 = 10.53271014302926
$ws.Cells.Item(15, 13).Value2 = 73.31390111996774

$ws.Cells.Item(16, 3).Value2 = 3.744823188419249
$ws.Cells.Item(16, 4).Value2 = 8.722097194435136
$ws.Cells.Item(16, 5).Value2 = 15.97959721953352
$ws.Cells.Item(16, 6).Value2 = 28.18689825325092
$ws.Cells.Item(16, 7).Value2 = 34.37525675426197
$ws.Cells.Item(16, 8).Value2 = 14.96694008110065
$ws.Cells.Item(16, 9).Value2 = 22.0765784131979
$ws.Cells.Item(16, 10).Value2 = 10.42769397000589
$ws.Cells.Item(16, 13).Value2 = 71.23195771146602

$ws.Cells.Item(17, 3).Value2 = 3.745775174537191
$ws.Cells.Item(17, 4).Value2 = 8.666636922476226
$ws.Cells.Item(17, 5).Value2 = 15.79849655355542
$ws.Cells.Item(17, 6).Value2 = 28.41601566832907
$ws.Cells.Item(17, 7).Value2 = 34.449954287468
$ws.Cells.Item(17, 8).Value2 = 15.0403076647286
$ws.Cells.Item(17, 9).Value2 = 22.15641237244655
$ws.Cells.Item(17, 10).Value2 = 10.36397714132209
$ws.Cells.Item(17, 13).Value2 = 69.92361505427385

$ws.Cells.Item(18, 3).Value2 = 3.746469495319036
$ws.Cells.Item(18, 4).Value2 = 8.634650139394195
$ws.Cells.Item(18, 5).Value2 = 15.69388773032675
$ws.Cells.Item(18, 6).Value2 = 28.55070981570385
$ws.Cells.Item(18, 7).Value2 = 34.49848091630102
$ws.Cells.Item(18, 8).Value2 = 15.08347239095901
$ws.Cells.Item(18, 9).Value2 = 22.20440031058893
$ws.Cells.Item(18, 10).Value2 = 10.32759667466845
$ws.Cells.Item(18, 13).Value2 = 69.15965313074592

$ws.Cells.Item(19, 3).Value2 = 3.746729605628381
$ws.Cells.Item(19, 4).Value2 = 8.623805758846007
$ws.Cells.Item(19, 5).Value2 = 15.65839576526673
$ws.Cells.Item(19, 6).Value2 = 28.59680787022344
$ws.Cells.Item(19, 7).Value2 = 34.51585292155796
$ws.Cells.Item(19, 8).Value2 = 15.09825149521892
$ws.Cells.Item(19, 9).Value2 = 22.22099925156579
$ws.Cells.Item(19, 10).Value2 = 10.31532561461797
$ws.Cells.Item(19, 13).Value2 = 68.89902517873843

$ws.Cells.Item(20, 3).Value2 = 3.745658601245682
$ws.Cells.Item(20, 4).Value2 = 8.672550013942386
$ws.Cells.Item(20, 5).Value2 = 15.81782176991756
$ws.Cells.Item(20, 6).Value2 = 28.39132263649354
$ws.Cells.Item(20, 7).Value2 = 34.44142366470096
$ws.Cells.Item(20, 8).Value2 = 15.03239724324604
$ws.Cells.Item(20, 9).Value2 = 22.14769867163672
$ws.Cells.Item(20, 10).Value2 = 10.3707323719723
$ws.Cells.Item(20, 13).Value2 = 70.06407432050506

$ws.Cells.Item(21, 3).Value2 = 3.743827270151343
$ws.Cells.Item(21, 4).Value2 = 8.835218934946036
$ws.Cells.Item(21, 5).Value2 = 16.34793949767442
$ws.Cells.Item(21, 6).Value2 = 27.73626502507879
$ws.Cells.Item(21, 7).Value2 = 34.25880775604897
$ws.Cells.Item(21, 8).Value2 = 14.82279326447054
$ws.Cells.Item(21, 9).Value2 = 21.92647446072735
$ws.Cells.Item(21, 10).Value2 = 10.56003551941077
$ws.Cells.Item(21, 13).Value2 = 73.84179049866233

$ws.Cells.Item(22, 3).Value2 = 3.743978788686041
$ws.Cells.Item(22, 4).Value2 = 8.940575374883625
$ws.Cells.Item(22, 5).Value2 = 16.68979693270017
$ws.Cells.Item(22, 6).Value2 = 27.3358033263752
$ws.Cells.Item(22, 7).Value2 = 34.1919680064177
$ws.Cells.Item(22, 8).Value2 = 14.69477439392018
$ws.Cells.Item(22, 9).Value2 = 21.80132377780873
$ws.Cells.Item(22, 10).Value2 = 10.68597451180626
$ws.Cells.Item(22, 13).Value2 = 76.20885922383106

$ws.Cells.Item(23, 3).Value2 = 3.743772500247115
$ws.Cells.Item(23, 4).Value2 = 8.884440581034855
$ws.Cells.Item(23, 5).Value2 = 16.50779171881645
$ws.Cells.Item(23, 6).Value2 = 27.54691825302253
$ws.Cells.Item(23, 7).Value2 = 34.22269682038805
$ws.Cells.Item(23, 8).Value2 = 14.76226112660069
$ws.Cells.Item(23, 9).Value2 = 21.86629180288943
$ws.Cells.Item(23, 10).Value2 = 10.61856342306964
$ws.Cells.Item(23, 13).Value2 = 74.95482845473963

$ws.Cells.Item(24, 3).Value2 = 3.745710846579929
$ws.Cells.Item(24, 4).Value2 = 8.669877020781287
$ws.Cells.Item(24, 5).Value2 = 15.80908636347386
$ws.Cells.Item(24, 6).Value2 = 28.40247711915426
$ws.Cells.Item(24, 7).Value2 = 34.44526300604506
$ws.Cells.Item(24, 8).Value2 = 15.03597048161823
$ws.Cells.Item(24, 9).Value2 = 22.15163163589706
$ws.Cells.Item(24, 10).Value2 = 10.36767754794465
$ws.Cells.Item(24, 13).Value2 = 70.00060937182494

$ws.Cells.Item(25, 3).Value2 = 3.753389604066787
$ws.Cells.Item(25, 4).Value2 = 8.435207232307592
$ws.Cells.Item(25, 5).Value2 = 15.03891367347362
$ws.Cells.Item(25, 6).Value2 = 29.43586036954019
$ws.Cells.Item(25, 7).Value2 = 34.8960965659233
$ws.Cells.Item(25, 8).Value2 = 15.36794130023386
$ws.Cells.Item(25, 9).Value2 = 22.53798091075377
$ws.Cells.Item(25, 10).Value2 = 10.10739439017206
$ws.Cells.Item(25, 13).Value2 = 64.21974629959919

Write-Output "Updated loading_percent values for rows 2-25"